$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10: group names shift to the next shared-string
# slot and report texts are replaced with the new report for that group.
$ws.Range('A2').Value = 'GRUPO DE INVESTIGACIONES AMBIENTALES'
$ws.Range('B2').Value = '207. Informe técnico: Manual para manejo de datos resultantes del monitoreo de calidad y cantidad de agua (Handeiling databank. Monitoringsgegevens waterkwaliteit en kwantiteit)  Año: 2003, Mes: 1, Idioma: Outros, Ciudad: Lelystad  Disponibilidad: Irrestricta, Número de páginas: 75, Número del contrato: , Institución en la cual prestó el servicio: Rijkswaterstaat 
 208. Informe técnico: INFORME DE CALIDAD DEL AGUA DEL RIO MEDELLIN DE LOS MONITOREOS 1, 2 Y 3  Año: 2004, Mes: 1, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restricta, Número de páginas: 90, Número del contrato: , Institución en la cual prestó el servicio: Area Metropolitana Del Valle De Aburrá'

$ws.Range('A3').Value = 'GISA'
$ws.Range('B3').Value = '1. Informe técnico: Realizar monitoreo de la calidad de aire (partículas menores de 10 micras (PM-10) sediméntales y gases SO2 y NO2); muestreo y análisis de la calidad del agua (cloruros, sulfatos, coliformes fecales, totales y mesofilicos), lluvia acida, (pH, acido sulfúrico H2SO4, acido nítrico HNO3 y acido clorhídr  Año: 2010, Mes: 1, Idioma: Español, Ciudad: HATONUEVO  Disponibilidad: No restringido, Número de páginas: 0, Número del contrato: Contarto 072 de 2010, Institución en la cual prestó el servicio: Corporación Autónoma Regional De La Guajira'

$ws.Range('A4').Value = 'Grupo de Investigación en calidad del agua y modelación hídrica (GICAMH)'
$ws.Range('B4').Value = '8. Informe técnico: FORTALECIMIENTO CRITERIOS ASOCIADOS A LA APROVISIONAMIENTO DEL AGUA/CALIDAD DEL AGUA PARA LA TOMA DE DECISIONES EN COLOMBIA. PROYECTO 2.1: DISENO DE UNA RED DE MONITOREO DE LA CALIDAD DEL AGUA SUPERFICIAL PARA EL AREA HIDROGRAFICA MAGDALENA-CAUCA  Año: 2016, Mes: 4, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restringido, Número de páginas: 55, Número del contrato: 81042877, Institución en la cual prestó el servicio: CORPORACIÓN CENTRO DE CIENCIA Y TECNOLOGÍA DE ANTIOQUIA 
 9. Informe técnico: FORTALECIMIENTO CRITERIOS ASOCIADOS A LA APROVISIONAMIENTO DEL AGUA/CALIDAD DEL AGUA PARA LA TOMA DE DECISIONES EN COLOMBIA. PROYECTO 2.2: HOJA DE RUTA PARA EL FORTALECIMIENTO DEL MONITOREO DE LA CALIDAD DEL AGUA SUPERFICIAL EN COLOMBIA  Año: 2016, Mes: 4, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restringido, Número de páginas: 55, Número del contrato: 81042877, Institución en la cual prestó el servicio: CORPORACIÓN CENTRO DE CIENCIA Y TECNOLOGÍA DE ANTIOQUIA'

$ws.Range('A5').Value = 'ECOLOGIA, INGENIERIA Y SOCIEDAD - EIS'
$ws.Range('B5').Value = '8. Informe técnico: Monitoreo quebrada Agua azul-111  Año: 2014, Mes: 8, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: No restringido, Número de páginas: 0, Número del contrato: 111-2014, Institución en la cual prestó el servicio: Compañía de Servicios Públicos Domiciliarios S.A E.S.P 
 65. Informe técnico: Monitoreo hidrobiológico y calidad del agua del río Manso-PMA 2012  Año: 2012, Mes: 1, Idioma: Español, Ciudad: MANIZALES  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: , Institución en la cual prestó el servicio:'

$ws.Range('A6').Value = 'Agua y Saneamiento'
$ws.Range('B6').Value = '3. Informe técnico: Diseño, implementación, evaluación, monitoreo y divulgación de sistemas eficientes de gestión del agua en hatos de producción lechera, de la subregión 1 del departamento de Risaralda  Año: 2007, Mes: 12, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Resolución Nº 1068 de 2007, Institución en la cual prestó el servicio: Gobernación De Risaralda 
 15. Informe técnico: Ejecución de pruebas de bombeo de larga duración en la red de monitoreo de agua subterráneas del sistema de acuíferos de Pereira  Año: 2011, Mes: 3, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio 481 de 2011, Institución en la cual prestó el servicio: Corporación Autónoma Regional De Risaralda - Carder 
 27. Informe técnico: Prestación de servicios de monitoreo para caracterización de aguas subterránea en 12 puntos para la Minera Seafield S.A.S  Año: 2013, Mes: 5, Idioma: Español, Ciudad: MANIZALES  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Carta de aceptacion, Institución en la cual prestó el servicio: Minera Seafield S.A.S 
 28. Informe técnico: Elaboración de estudio de caracterización de Aguas de Pozos de Monitoreo del Ingenio Risaralda S.A.  Año: 2013, Mes: 2, Idioma: Español, Ciudad: LA VIRGINIA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Carta de aceptacion, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 43. Informe técnico: Implementación del protocolo IDEAM de monitoreo y seguimiento de aguas continentales en el departamento de Caldas, como insumo para la ordenación del recurso hídrico  Año: 2008, Mes: 11, Idioma: Español, Ciudad: MANIZALES  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Contrato 158 de 2008, Institución en la cual prestó el servicio: Corporación Autónoma Regional Para El Desarrollo De Caldas - Corpocaldas (Antes Cramsa) 
 49. Informe técnico: Formulacion del plan de manejo del acuifero de la de la cuenca del rio Risaralda, implementacion de los planes de manejo de los acuiferos Santagueda km 41- Irra y rio gande la Magdalena y operación de la red de monitoreo de agua subterranea de los acuiferos en la jurisdiccion de Caldas.  Año: 2015, Mes: 8, Idioma: Español, Ciudad: MANIZALES  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Contrato 136 de 2015, Institución en la cual prestó el servicio: Corporación Autónoma Regional Para El Desarrollo De Caldas - Corpocaldas (Antes Cramsa) 
 53. Informe técnico: Capacitación, monitoreo y seguimiento de los procesos de dragado en los embalses Cameguadua, San Francisco, plantas menores y sistemas de tratamiento de aguas residuales de la CHEC S.A. E.S.P.  Año: 2012, Mes: 6, Idioma: Español, Ciudad: MANIZALES  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Contrato 2267 de 2012, Institución en la cual prestó el servicio: Corporación Autónoma Regional Para El Desarrollo De Caldas - Corpocaldas (Antes Cramsa) 
 63. Informe técnico: Estudio de caracterizacion de aguas de pozos de monitoreo  Año: 2014, Mes: 7, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Cotizacion, Institución en la cual prestó el servicio: 
 73. Informe técnico: Informe estudio de caracterizacion de agua de pozos de monitoreo semestre I año 2015 Ingenio Risaralda S.A.  Año: 2015, Mes: 7, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 012-15- de 2015, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 74. Informe técnico: Informe estudio de caracterizacion de agua de pozos de monitoreo semestre II a¿o 2015 Ingenio Risaralda S.A.  Año: 2015, Mes: 12, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 097-15- de 2015, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 79. Informe técnico: Estudios de caracterización de aguas de pozos de monitoreo año 2016  Año: 2016, Mes: 5, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Propuesta t¿cnico econ¿mica 018-16, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 81. Informe técnico: Estudio de caracterizacion de agua de pozos de monitoreo en el Ingenio Risaralda S.A, segundo semestre de 2016  Año: 2016, Mes: 9, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Propuesta t?cnico econ?mica 106-16, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 100. Informe técnico: Estudio de caracterización de agua de pozos de monitoreo semestre I año 2018 Ingenio Risaralda S.A.  Año: 2018, Mes: 6, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 58, Número del contrato: Orden de compra # 3620172063, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 106. Informe técnico: Estudio de caracterización de agua de pozos de monitoreo semestre II año 2018 Ingenio Risaralda S.A.  Año: 2018, Mes: 12, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 40, Número del contrato: Orden de compra # 3620178418, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 127. Informe técnico: Caracterización de agua subterránea para once (11) pozos de monitoreo ubicados en el Departamento de Risaralda. 2019-I  Año: 2019, Mes: 5, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 024-19, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 133. Informe técnico: Caracterización de agua subterránea para once (11) pozos de monitoreo ubicados en el Departamento de Risaralda. 2019-II  Año: 2019, Mes: 11, Idioma: Español, Ciudad: PEREIRA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 103-19, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 137. Informe técnico: Caracterización de agua subterránea proveniente de 11 pozos de monitoreo ubicados en el Departamento de Risaralda 2020  Año: 2020, Mes: 10, Idioma: Español, Ciudad: LA VIRGINIA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 022-20, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A. 
 161. Informe técnico: Estudio de caracterizacion de aguas de pozos de monitoreo  Año: 2014, Mes: 7, Idioma: Español, Ciudad: LA VIRGINIA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Carta de aceptacion de propuesta, Institución en la cual prestó el servicio: INGENIO RISARALDA S.A.'

$ws.Range('A7').Value = 'INVEMAR - Calidad Ambiental Marina'
$ws.Range('B7').Value = '1. Informe técnico: Fortalecimiento al monitoreo calidad del agua, del ecosistema de manglar y de los recursos pesqueros de la Ciénaga Grande de Santa Marta  Año: 2015, Mes: 1, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio 040/14, Institución en la cual prestó el servicio: 
 10. Informe técnico: Implementar un proyecto piloto sobre el protocolo de monitoreo y seguimiento del agua marino-costera.  Año: 2007, Mes: 7, Idioma: Español, Ciudad: BOGOTÁ, D.C.  Disponibilidad: Restringido, Número de páginas: 200, Número del contrato: Contrato No C-0550-06. Convenio No 001/04, Institución en la cual prestó el servicio: 
 34. Informe técnico: Informe técnico final del Monitoreo de la calidad de aguas y sedimentos, y de los ecosistemas marinos en las zonas afectadas por el derrame de aceite de palma, bahías de Taganga y Santa Marta (Magdalena). Realizado para la firma TERLICA S.A  Año: 2008, Mes: 1, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: Restringido, Número de páginas: 82, Número del contrato: , Institución en la cual prestó el servicio: TERLICA 
 35. Informe técnico: Monitoreo de calidad de aguas, sedimentos y aspectos biológicos en una zona industrial de Mamonal, Bahía de Cartagena  Año: 2008, Mes: 1, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: Restringido, Número de páginas: 189, Número del contrato: , Institución en la cual prestó el servicio: Aribasa 
 43. Informe técnico: Monitoreo ambiental de la calidad de aguas, sedimentos y comunidades marinas en la zona de influencia directa del pozo exploratorio Arazá I  Año: 2007, Mes: 6, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: Restringido, Número de páginas: 272, Número del contrato: , Institución en la cual prestó el servicio: 
 48. Informe técnico: Monitoreo de calidad de aguas, sedimentos y comunidades biológicas durante la perforación de pozos en Chuchupa B  Año: 2006, Mes: 1, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: Restricta, Número de páginas: 120, Número del contrato: , Institución en la cual prestó el servicio: Chevron Texaco 
 56. Informe técnico: Fortalecimiento al monitoreo calidad del agua, del ecosistema de manglar y de los recursos pesqueros de la Ciénaga Grande de Santa Marta. Informe Técnico Final. Santa Marta.  Año: 2015, Mes: 12, Idioma: Español, Ciudad: Disponibilidad: No restringido, Número de páginas: 58, Número del contrato: Convenio 040/14, Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR 
 62. Informe técnico: FORTALECIMIENTO DEL MONITOREO DE LA CALIDAD DEL AGUA, DEL ECOSISTEMA DE MANGLAR Y DE LOS RECURSOS PESQUEROS DE LA CIENAGA GRANDE DE SANTA MARTA Año: 2015, Mes: 1, Idioma: Español, Ciudad: SANTA MARTA  Disponibilidad: No restringido, Número de páginas: 45, Número del contrato: , Institución en la cual prestó el servicio: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR 
 66. Informe técnico: Fortalecimiento al monitoreo calidad del agua, del ecosistema de manglar y de los recursos pesqueros de la Ci¿naga Grande de Santa Marta  Año: 2014, Mes: 1, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: , Institución en la cual prestó el servicio: 
 69. Informe técnico: Fortalecimiento al monitoreo calidad del agua, del ecosistema de manglar y de los recursos pesqueros de la Ci¿naga Grande de Santa Marta  Año: 2015, Mes: 1, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: , Institución en la cual prestó el servicio:'

$ws.Range('A8').Value = 'Grupo de Investigación en Modelación y Evaluación de Sistemas Ambientales. GiMESA'
$ws.Range('B8').Value = '7. Informe técnico: Informe No.2. Propuesta de red de monitoreo de aguas atmosféricas, superficiales y subterráneas y medición de parámetrso de campo del estudio de impacto ambiental del proyecto de aprovechamiento y recursos minerales en los cerros Queresas y Porvenir  Año: 2019, Mes: 2, Idioma: Español, Ciudad: MONTERÍA  Disponibilidad: Restringido, Número de páginas: 38, Número del contrato: 4541319879, Institución en la cual prestó el servicio: Cerro Matoso S.A 
 26. Informe técnico: Informe N° 2. Propuesta de Red de Monitoreo de Aguas Atmosféricas, Superficiales y Subterráneas y Medición de Parámetros de Campo del Estudio de Impacto Ambiental del Proyecto Aprovechamiento de Recursos Minerales en los Cerros Queresas y Porvenir  Año: 2019, Mes: 2, Idioma: Español, Ciudad: Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: 4541319879, Institución en la cual prestó el servicio: Cerro Matoso S.A 
 39. Informe técnico: Informe N° 2. Propuesta de Red de Monitoreo de Aguas Atmosféricas, Superficiales y Subterráneas y Medición de Parámetros de Campo del Estudio de Impacto Ambiental del Proyecto Aprovechamiento de Recursos Minerales en los Cerros Queresas y Porvenir  Año: 2019, Mes: 2, Idioma: Español, Ciudad: MONTERÍA  Disponibilidad: Restringido, Número de páginas: 38, Número del contrato: 4541319879, Institución en la cual prestó el servicio: Cerro Matoso S.A'

$ws.Range('A9').Value = 'Gestión Integrada de Recursos Hídricos'
$ws.Range('B9').Value = '11. Informe técnico: Monitoreo de las Alianzas para el Aprendizaje del proyecto usos múltiples del agua  Año: 2006, Mes: 1, Idioma: Español, Ciudad: CALI  Disponibilidad: Irrestricta, Número de páginas: 41, Número del contrato: , Institución en la cual prestó el servicio: Challenge Program On Water And Food 
 29. Informe técnico: Monitoreo de las Alianzas para el Aprendizaje del proyecto Usos Múltiples del Agua  Año: 2006, Mes: 1, Idioma: Español, Ciudad: Disponibilidad: No restringido, Número de páginas: 41, Número del contrato: , Institución en la cual prestó el servicio:'

$ws.Range('A10').Value = 'Grupo de Estudios Ambientales para la Sostenibilidad, la Innovación y el Desarrollo GEASID'
$ws.Range('B10').Value = '33. Informe técnico: Determinación de los índices de calidad y contaminación de agua ICA''s e ICO''s de la primera campaña de monitoreo de las quebradas la Gómez, Santos Gutiérrez, Pescado, Islitas, Caño Peruétano y Ciénaga de Paredes, en los municipio de Sabana de Torres y Puerto Wilches, Santander.  Año: 2015, Mes: 6, Idioma: Español, Ciudad: SAN GIL  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio # 5219582, Institución en la cual prestó el servicio: ECOPETROL S. A. 
 50. Informe técnico: Determinación de los índices de calidad y contaminación de agua ICA''s e ICO''s de la tercera campaña de monitoreo de las quebradas la Gómez, Santos Gutiérrez, Pescado, Islitas, Caño Peruétano y Ciénaga de Paredes, en los municipio de Sabana de Torres y Puerto Wilches, Santander.  Año: 2016, Mes: 1, Idioma: Español, Ciudad: SAN GIL  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio # 5219582, Institución en la cual prestó el servicio: ECOPETROL S. A. 
 75. Informe técnico: Determinación de los índices de calidad y contaminación de agua ICA''s e ICO''s de la primera campaña de monitoreo de las quebradas la Gómez, Santos Gutiérrez, Pescado, Islitas, Caño Peruétano y Ciénaga de Paredes, en los municipio de Sabana de Torres y Puerto Wilches, Santander.  Año: 2015, Mes: 6, Idioma: Español, Ciudad: SAN GIL  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio # 5219582, Institución en la cual prestó el servicio: ECOPETROL S. A. 
 77. Informe técnico: Determinación de los índices de calidad y contaminación de agua ICA''s e ICO''s de la tercera campaña de monitoreo de las quebradas la Gómez, Santos Gutiérrez, Pescado, Islitas, Caño Peruétano y Ciénaga de Paredes, en los municipio de Sabana de Torres y Puerto Wilches, Santander.  Año: 2016, Mes: 1, Idioma: Español, Ciudad: SAN GIL  Disponibilidad: Restringido, Número de páginas: 0, Número del contrato: Convenio # 5219582, Institución en la cual prestó el servicio: ECOPETROL S. A.'

# Add new row 11 for the new group, copying formatting from row 10's A cell
$ws.Range('A10').Copy()
$ws.Range('A11').PasteSpecial(-4122)
$ws.Range('A11').Value = 'Grupo de investigación en Agua y Medio Ambiente - AMA'
$ws.Range('B11').Value = '78. Informe técnico: FORTALECIMIENTO CRITERIOS ASOCIADOS A LA APROVISIONAMIENTO DEL AGUA/CALIDAD DEL AGUA PARA LA TOMA DE DECISIONES EN COLOMBIA. PROYECTO 2.1: DISEÑO DE UNA RED DE MONITOREO DE LA CALIDAD DEL AGUA SUPERFICIAL PARA EL AREA HIDROGRAFICA MAGDALENA - CAUCA  Año: 2016, Mes: 4, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restringido, Número de páginas: 55, Número del contrato: 81042877, Institución en la cual prestó el servicio: CORPORACIÓN CENTRO DE CIENCIA Y TECNOLOGÍA DE ANTIOQUIA 
 79. Informe técnico: FORTALECIMIENTO CRITERIOS ASOCIADOS A LA APROVISIONAMIENTO DEL AGUA/CALIDAD DEL AGUA PARA LA TOMA DE DECISIONES EN COLOMBIA. PROYECTO 2.2: HOJA DE RUTA PARA EL FORTALECIMIENTO DEL MONITOREO DE LA CALIDAD DEL AGUA SUPERFICIAL EN COLOMBIA  Año: 2016, Mes: 4, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restringido, Número de páginas: 55, Número del contrato: 81042877, Institución en la cual prestó el servicio: CORPORACIÓN CENTRO DE CIENCIA Y TECNOLOGÍA DE ANTIOQUIA 
 94. Informe técnico: FORTALECIMIENTO CRITERIOS ASOCIADOS A LA APROVISIONAMIENTO DEL AGUA/CALIDAD DEL AGUA PARA LA TOMA DE DECISIONES EN COLOMBIA. PROYECTO 2.2: HOJA DE RUTA PARA EL FORTALECIMIENTO DEL MONITOREO DE LA CALIDAD DEL AGUA SUPERFICIAL EN COLOMBIA  Año: 2016, Mes: 4, Idioma: Español, Ciudad: MEDELLÍN  Disponibilidad: Restringido, Número de páginas: 55, Número del contrato: 81042877, Institución en la cual prestó el servicio: CORPORACIÓN CENTRO DE CIENCIA Y TECNOLOGÍA DE ANTIOQUIA'